# Order-of-magnitude correction for 'time to scan' calculations. Divide by 1000.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# EO FOV / Horizontal / Vertical resolution inputs
$ws.Range("C13").Value = 60
$ws.Range("C14").Value = 1800
$ws.Range("C15").Value = 1800

# IR FOV / Horizontal / Vertical resolution inputs
$ws.Range("C20").Value = 90
$ws.Range("C21").Value = 1200
$ws.Range("C22").Value = 1200

# Time to scan formulas: divide the numerator by 1000 (25000000 -> 25000)
$ws.Range("F17").Formula = "=25000*60/EO_Ground_Coverage_Rate"
$ws.Range("F33").Formula = "=25000*60/IR_Ground_Coverage_Rate"

# Update the active selection to match the recorded sheet view
$ws.Range("C16").Select()
